$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("add_device_to_group_failed", "Ocurrió un error al agregar el dispositivo al grupo de dispositivos", "An error occurred while adding the device to the device group"),
    @("get_devices_by_owner_failed", "Ocurrió un error al recuperar los dispositivos del usuario", "An error occurred while retrieving user devices"),
    @("remove_device_from_group_failed", "Ocurrió un error al eliminar el dispositivo del grupo de dispostivos.", "An error occurred while removing the device from the device group."),
    @("save_device_failed", "Ocurrió un error al guardar la información del dispositivo.", "An error occurred while saving device information."),
    @("update_device_failed", "Ocurrió un error al actualizar la información del dispositivo.", "An error occurred while updating device information."),
    @("device_added_successfully", "La solicitud de agregación del dispositivo al grupo de dispositivos del usuario ha sido registrada con éxito.", "The request to add the device to the user's device group has been successfully registered."),
    @("device_updated_successfully", "La solicitud de actualización del dispositivo ha sido registrada con éxito.", "The device update request has been successfully registered."),
    @("device_saved_successfully", "La solicitud de  almacenamiento del dispositivo ha sido registrada con éxito.", "The device storage request has been registered successfully."),
    @("device_deleted_successfully", "La solicitud de borrado del dispositivo ha sido registrada con éxito.", "The device wipe request has been successfully registered.")
)

$startRow = 104

# Shared strings must land in the workbook's string table in the same
# order the source commit produced them: all the NAME keys (column A)
# first, then the ES/EN message pairs (columns B/C) row by row.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
}

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$ws.Range("C112").Select()
